$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.004.09'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').Value = '1.621.28'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.25'
$ws.Range('E5').Value = '  -1.61%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('E9').Value = '  -1.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.90'
$ws.Range('E10').Value = '  -0.76%  '
$ws.Range('E11').Value = '  -1.38%  '
$ws.Range('D12').Value = '1.847.55'
$ws.Range('E12').Value = '  -1.00%  '
$ws.Range('D13').Value = '1.619.07'
$ws.Range('E13').Value = '  -1.04%  '
$ws.Range('E14').Value = '  -0.59%  '
$ws.Range('E15').Value = '  -0.69%  '
$ws.Range('D16').Value = '26.988.60'
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.26'
$ws.Range('E17').Value = '  -3.41%  '
$ws.Range('D18').Value = '0.0₃0735'
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '213.73'
$ws.Range('E19').Value = '  -1.27%  '
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('E22').Value = '  -2.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.35'
$ws.Range('E23').Value = '  -7.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.94'
$ws.Range('E24').Value = '  -1.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.72'
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.46'
$ws.Range('E26').Value = '  +1.04%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('E28').Value = '  -3.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.49'
$ws.Range('E29').Value = '  -1.00%  '
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('E31').Value = '  -1.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.28'
$ws.Range('E32').Value = '  -2.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.699'
$ws.Range('E33').Value = '  +26.88%  '
$ws.Range('E34').Value = '  -1.08%  '
$ws.Range('D35').Value = '1.334.46'
$ws.Range('E35').Value = '  +2.49%  '
$ws.Range('E36').Value = '  -0.81%  '
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('E38').Value = '  -0.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.839'
$ws.Range('E39').Value = '  -1.65%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('E41').Value = '  -1.20%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.21'
$ws.Range('E42').Value = '  -1.31%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.34'
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.74'
$ws.Range('E44').Value = '  +2.61%  '
$ws.Range('D45').Value = '1.758.84'
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.79'
$ws.Range('E46').Value = '  -1.64%  '
$ws.Range('E47').Value = '  +2.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.843'
$ws.Range('E48').Value = '  +16.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0516'
$ws.Range('E49').Value = '  +0.27%  '
$ws.Range('E51').Value = '  -0.81%  '
